# fix: revert admin dev default; seed customers only when table empty;
#      autosave on customer select when hours/day present
#
# Applies to Chris_Zavesky_2025-12-29.xlsx:
#   1. Rewrites the 4 timesheet detail rows on "Weekly Timesheet" with the
#      corrected (reverted) data set + recomputed subtotal/grand-total block.
#   2. Adds a "Jason Schema" sheet holding the same rows in a flattened,
#      per-employee schema (Employee / Employee ID / Date / Client / Hours /
#      Rate / Total / Type / Notes) used by the seeding job.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Weekly Timesheet"

# Row 4 used to be the (bold/orange) SUBTOTAL row and row 5 the "Category"
# marker row; both are being reused as plain detail / new summary rows
# below, so strip any leftover formatting (including the row-level default
# style) from the whole block first and rebuild only the formatting that is
# actually wanted for the new layout.
for ($r = 2; $r -le 12; $r++) {
    $ws.Rows.Item($r).ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. Weekly Timesheet - detail rows (2-5)
# ---------------------------------------------------------------------------

$data = @(
    @("2025-12-30", "McGill", 7,  "Regular"),
    @("2025-12-31", "Hall",   8,  "Regular"),
    @("2026-01-01", "Tormey", 20, "Holiday"),
    @("2026-01-02", "Richer", 20, "Regular")
)

$r = 2
foreach ($row in $data) {
    # Leading apostrophe forces these date-shaped strings to stay literal
    # text (matching the source data's shared-string "General" cells)
    # instead of being auto-parsed into Excel date serials.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $r++
}

# row 6 stays blank (spacer row) - deliberately left untouched

# ---------------------------------------------------------------------------
# 2. SUBTOTAL block (row 7) + "Category: ADMIN" marker (row 8)
# ---------------------------------------------------------------------------

$totalHours = 0
foreach ($row in $data) { $totalHours += $row[2] }

$ws.Cells.Item(7, 1).Value = "SUBTOTAL"
$ws.Cells.Item(7, 2).Value = ""
$ws.Cells.Item(7, 3).Value = $totalHours
$ws.Cells.Item(7, 4).Value = "Reg: $totalHours / OT: 0"
$ws.Cells.Item(7, 5).Value = ""
$ws.Cells.Item(7, 6).Value = 0

$ws.Cells.Item(8, 1).Value = "Category: ADMIN"
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(8, 3).Value = ""
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(8, 6).Value = ""

# row 9 stays blank (spacer row) - deliberately left untouched

# ---------------------------------------------------------------------------
# 3. Grand-total block (rows 10-12)
# ---------------------------------------------------------------------------

$ws.Cells.Item(10, 1).Value = "HOURLY SUBTOTAL"
$ws.Cells.Item(10, 2).Value = ""
$ws.Cells.Item(10, 3).Value = ""
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = ""
$ws.Cells.Item(10, 6).Value = 0

$ws.Cells.Item(11, 1).Value = "ADMIN SUBTOTAL"
$ws.Cells.Item(11, 2).Value = ""
$ws.Cells.Item(11, 3).Value = ""
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = 0

$ws.Cells.Item(12, 1).Value = "GRAND TOTAL"
$ws.Cells.Item(12, 2).Value = ""
$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(12, 6).Value = 0

# ---------------------------------------------------------------------------
# 4. Formatting for the new rows
# ---------------------------------------------------------------------------

# SUBTOTAL row keeps the existing "orange" subtotal style (same as before).
$subtotalRange = $ws.Range("A7:F7")
$subtotalRange.Font.Bold = $true
$subtotalRange.Interior.Color = 11591935   # RGB(FFE0B0) in BGR long form

# New "hourly" / "admin" subtotal rows use a pale-gold fill.
$hourlyRange = $ws.Range("A10:F11")
$hourlyRange.Font.Bold = $true
$hourlyRange.Interior.Color = 14742522     # RGB(FAF3E0) in BGR long form

# Grand-total row uses a pale-green fill with bold red text.
$grandRange = $ws.Range("A12:F12")
$grandRange.Font.Bold = $true
$grandRange.Font.Color = 255               # RGB(FF0000) in BGR long form
$grandRange.Interior.Color = 14743784      # RGB(E8F8E0) in BGR long form

# Currency number format on the Rate/Total columns for the new rows.
# (rows 6 and 9 are spacer rows and are deliberately skipped so no empty
# <row> element gets materialised for them)
$ws.Range("E2:F5").NumberFormat = """$""#,##0.00"
$ws.Range("E7:F8").NumberFormat = """$""#,##0.00"
$ws.Range("E10:F12").NumberFormat = """$""#,##0.00"

# ---------------------------------------------------------------------------
# 5. "Jason Schema" worksheet - flattened per-employee seed table
# ---------------------------------------------------------------------------

$schema = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$schema.Name = "Jason Schema"

$headers = @("Employee", "Employee ID", "Date", "Client", "Hours", "Rate", "Total", "Type", "Notes")
for ($c = 1; $c -le $headers.Length; $c++) {
    $schema.Cells.Item(1, $c).Value = $headers[$c - 1]
}
$schema.Range("A1:I1").Font.Bold = $true
# Rate/Total headers carry the currency format too (column default).
$schema.Range("F1:G1").NumberFormat = """$""#,##0.00"

$employee = "Chris Zavesky"
$employeeId = "emp_5chpvt65"

$r = 2
foreach ($row in $data) {
    $schema.Cells.Item($r, 1).Value = $employee
    $schema.Cells.Item($r, 2).Value = $employeeId
    $schema.Cells.Item($r, 3).Value = "'" + $row[0]
    $schema.Cells.Item($r, 4).Value = $row[1]
    $schema.Cells.Item($r, 5).Value = $row[2]
    $schema.Cells.Item($r, 6).Value = 0
    $schema.Cells.Item($r, 7).Value = 0
    $schema.Cells.Item($r, 8).Value = $row[3]
    $schema.Cells.Item($r, 9).Value = ""
    $r++
}

$schema.Range("F2:G5").NumberFormat = """$""#,##0.00"

$schema.Columns.Item(1).ColumnWidth = 20 - (5/6)
$schema.Columns.Item(2).ColumnWidth = 18 - (5/6)
$schema.Columns.Item(3).ColumnWidth = 12 - (5/6)
$schema.Columns.Item(4).ColumnWidth = 25 - (5/6)
$schema.Columns.Item(5).ColumnWidth = 8 - (5/6)
$schema.Columns.Item(6).ColumnWidth = 10 - (5/6)
$schema.Columns.Item(7).ColumnWidth = 12 - (5/6)
$schema.Columns.Item(8).ColumnWidth = 10 - (5/6)
$schema.Columns.Item(9).ColumnWidth = 30 - (5/6)

$ws.Activate()
$ws.Range("A1").Select()
